$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Adapt the "funcionarios" template to receive the exported rows (3 new
# employees appended below the existing data: Pedro, Charles, Deb).
# ---------------------------------------------------------------------------

# --- Row 4: Pedro ----------------------------------------------------------
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = "Pedro"
$ws.Cells.Item(4, 4).Value = "Dev front-end"
$ws.Cells.Item(4, 5).Value = "'5000.00"
$ws.Cells.Item(4, 6).Value = 45664
$ws.Cells.Item(4, 6).NumberFormat = $ws.Cells.Item(2, 6).NumberFormat
$ws.Cells.Item(4, 7).Value = 2
$ws.Cells.Item(4, 8).Value = "carta dev front "

# --- Row 5: Charles ---------------------------------------------------------
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(5, 3).Value = "Charles"
$ws.Cells.Item(5, 4).Value = "Dev front-end"
$ws.Cells.Item(5, 5).Value = "'10000.00"
$ws.Cells.Item(5, 6).Value = 45575
$ws.Cells.Item(5, 6).NumberFormat = $ws.Cells.Item(2, 6).NumberFormat
$ws.Cells.Item(5, 7).Value = 3
$ws.Cells.Item(5, 8).Value = ""

# --- Row 6: Deb --------------------------------------------------------------
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 3
$ws.Cells.Item(6, 3).Value = "Deb"
$ws.Cells.Item(6, 4).Value = "Dev front-end"
$ws.Cells.Item(6, 5).Value = "'5000.00"
$ws.Cells.Item(6, 6).Value = 45209
$ws.Cells.Item(6, 6).NumberFormat = $ws.Cells.Item(2, 6).NumberFormat
$ws.Cells.Item(6, 7).Value = 2
$ws.Cells.Item(6, 8).Value = ""
